# Auto-generated script applying 2024-03-17 crime data update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 7705
$ws.Cells.Item(2, 11).Value = 1466
$ws.Cells.Item(3, 10).Value = 8078
$ws.Cells.Item(3, 11).Value = 1388
$ws.Cells.Item(4, 10).Value = 1794
$ws.Cells.Item(4, 11).Value = 304
$ws.Cells.Item(5, 11).Value = 94
$ws.Cells.Item(6, 11).Value = 1785
$ws.Cells.Item(7, 10).Value = 29262
$ws.Cells.Item(7, 11).Value = 5037

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(2, 11).Value = 17
$ws.Cells.Item(6, 11).Value = 44
$ws.Cells.Item(7, 11).Value = 78

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 11).Value = 93
$ws.Cells.Item(3, 11).Value = 89
$ws.Cells.Item(4, 11).Value = 17
$ws.Cells.Item(6, 11).Value = 99
$ws.Cells.Item(7, 11).Value = 305

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 10).Value = 176
$ws.Cells.Item(4, 11).Value = 5
$ws.Cells.Item(6, 10).Value = 162
$ws.Cells.Item(7, 11).Value = 100

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 11).Value = 59
$ws.Cells.Item(3, 11).Value = 80
$ws.Cells.Item(6, 11).Value = 51
$ws.Cells.Item(7, 11).Value = 204

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 11).Value = 30
$ws.Cells.Item(7, 11).Value = 86

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 11).Value = 35
$ws.Cells.Item(3, 11).Value = 55
$ws.Cells.Item(6, 11).Value = 63
$ws.Cells.Item(7, 11).Value = 171

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 11).Value = 29
$ws.Cells.Item(7, 11).Value = 92

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(6, 11).Value = 7
$ws.Cells.Item(7, 11).Value = 19

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 11).Value = 34
$ws.Cells.Item(7, 11).Value = 144
$ws.Cells.Item(8, 11).Value = 305
$ws.Cells.Item(9, 11).Value = 24
$ws.Cells.Item(10, 11).Value = 33
$ws.Cells.Item(11, 11).Value = 101
$ws.Cells.Item(14, 11).Value = 32
$ws.Cells.Item(15, 11).Value = 43
$ws.Cells.Item(20, 11).Value = 116
$ws.Cells.Item(23, 11).Value = 50
$ws.Cells.Item(27, 11).Value = 60
$ws.Cells.Item(29, 10).Value = 1556
$ws.Cells.Item(29, 11).Value = 229
$ws.Cells.Item(30, 11).Value = 19
$ws.Cells.Item(31, 11).Value = 57
$ws.Cells.Item(33, 11).Value = 204
$ws.Cells.Item(36, 11).Value = 56
$ws.Cells.Item(37, 11).Value = 171
$ws.Cells.Item(39, 11).Value = 7
$ws.Cells.Item(41, 11).Value = 55
$ws.Cells.Item(42, 11).Value = 175
$ws.Cells.Item(43, 11).Value = 48
$ws.Cells.Item(44, 11).Value = 45
$ws.Cells.Item(48, 11).Value = 55
$ws.Cells.Item(52, 11).Value = 137
$ws.Cells.Item(53, 11).Value = 78
$ws.Cells.Item(54, 11).Value = 87
$ws.Cells.Item(55, 11).Value = 52
$ws.Cells.Item(60, 11).Value = 39
$ws.Cells.Item(63, 11).Value = 17
$ws.Cells.Item(67, 11).Value = 194
$ws.Cells.Item(73, 11).Value = 50
$ws.Cells.Item(76, 10).Value = 415
$ws.Cells.Item(76, 11).Value = 68
$ws.Cells.Item(78, 11).Value = 72
$ws.Cells.Item(83, 11).Value = 100
$ws.Cells.Item(85, 11).Value = 258
$ws.Cells.Item(87, 10).Value = 98
$ws.Cells.Item(88, 11).Value = 65
$ws.Cells.Item(89, 11).Value = 65
$ws.Cells.Item(90, 11).Value = 50
$ws.Cells.Item(92, 11).Value = 22
$ws.Cells.Item(95, 11).Value = 86
$ws.Cells.Item(96, 11).Value = 66
$ws.Cells.Item(99, 11).Value = 92
$ws.Cells.Item(101, 10).Value = 29262
$ws.Cells.Item(101, 11).Value = 5037

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(6, 11).Value = 23
$ws.Cells.Item(7, 11).Value = 57

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 11).Value = 60
$ws.Cells.Item(7, 11).Value = 194

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(6, 11).Value = 33
$ws.Cells.Item(7, 11).Value = 87

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 11).Value = 59
$ws.Cells.Item(3, 11).Value = 74
$ws.Cells.Item(6, 10).Value = 398
$ws.Cells.Item(7, 10).Value = 1556
$ws.Cells.Item(7, 11).Value = 229

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(3, 11).Value = 11
$ws.Cells.Item(7, 11).Value = 55

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 11).Value = 18
$ws.Cells.Item(7, 11).Value = 45

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(3, 11).Value = 13
$ws.Cells.Item(4, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 415
$ws.Cells.Item(7, 11).Value = 68

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(2, 11).Value = 14
$ws.Cells.Item(7, 11).Value = 32

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(2, 11).Value = 16
$ws.Cells.Item(6, 11).Value = 26
$ws.Cells.Item(7, 11).Value = 55

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 11).Value = 76
$ws.Cells.Item(7, 11).Value = 175

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(6, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 33

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(6, 11).Value = 24
$ws.Cells.Item(7, 11).Value = 72

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 11).Value = 21
$ws.Cells.Item(6, 11).Value = 18
$ws.Cells.Item(7, 11).Value = 52

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(3, 11).Value = 15
$ws.Cells.Item(7, 11).Value = 50

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(3, 11).Value = 8
$ws.Cells.Item(7, 11).Value = 66

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 11).Value = 34
$ws.Cells.Item(7, 11).Value = 116

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(6, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 56

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 11).Value = 55
$ws.Cells.Item(3, 11).Value = 45
$ws.Cells.Item(7, 11).Value = 144

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 11).Value = 12
$ws.Cells.Item(7, 11).Value = 43

$ws = $wb.Worksheets.Item('Greektown')
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(6, 11).Value = 7

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(6, 11).Value = 43
$ws.Cells.Item(7, 11).Value = 101

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(7, 11).Value = 24

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(6, 11).Value = 22
$ws.Cells.Item(7, 11).Value = 50

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 11).Value = 13
$ws.Cells.Item(7, 11).Value = 34

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Cells.Item(3, 11).Value = 4
$ws.Cells.Item(7, 11).Value = 22

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 11).Value = 15
$ws.Cells.Item(7, 11).Value = 65

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(3, 11).Value = 24
$ws.Cells.Item(7, 11).Value = 65

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Cells.Item(4, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 60

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 11).Value = 23
$ws.Cells.Item(3, 11).Value = 11
$ws.Cells.Item(7, 11).Value = 50

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Cells.Item(3, 11).Value = 17
$ws.Cells.Item(6, 11).Value = 12
$ws.Cells.Item(7, 11).Value = 39

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 11).Value = 20
$ws.Cells.Item(7, 11).Value = 48

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 11).Value = 96
$ws.Cells.Item(3, 11).Value = 80
$ws.Cells.Item(7, 11).Value = 258

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(6, 11).Value = 67
$ws.Cells.Item(7, 11).Value = 137

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(3, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 98
